$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.04532070457935333
$ws.Range("B2").Value = 0.9880834817886353
$ws.Range("C2").Value = 0.01455355249345303
$ws.Range("D2").Value = 0.9979322552680969
$ws.Range("A3").Value = 0.006729683838784695
$ws.Range("B3").Value = 0.998579740524292
$ws.Range("C3").Value = 0.00825185514986515
$ws.Range("D3").Value = 0.9982503652572632
$ws.Range("A4").Value = 0.002779328729957342
$ws.Range("B4").Value = 0.9992205500602722
$ws.Range("C4").Value = 0.004820717498660088
$ws.Range("D4").Value = 0.9990456700325012
$ws.Range("A5").Value = 0.001849205465987325
$ws.Range("B5").Value = 0.9993937611579895
$ws.Range("C5").Value = 0.003767602378502488
$ws.Range("D5").Value = 0.9995228052139282
$ws.Range("A6").Value = 0.001067447359673679
$ws.Range("B6").Value = 0.9996882081031799
$ws.Range("C6").Value = 0.001933840336278081
$ws.Range("D6").Value = 0.9996818900108337
$ws.Range("A7").Value = 0.0009955903515219688
$ws.Range("B7").Value = 0.9997401833534241
$ws.Range("C7").Value = 0.002882387489080429
$ws.Range("D7").Value = 0.9995228052139282
$ws.Range("A8").Value = 0.0005075388471595943
$ws.Range("B8").Value = 0.9999133944511414
$ws.Range("C8").Value = 0.002420859644189477
$ws.Range("D8").Value = 0.9998409152030945
$ws.Range("A9").Value = 0.0009402516880072653
$ws.Range("B9").Value = 0.9997748136520386
$ws.Range("C9").Value = 0.003609405132010579
$ws.Range("D9").Value = 0.9995228052139282
$ws.Range("A10").Value = 0.0006790620973333716
$ws.Range("B10").Value = 0.9998267889022827
$ws.Range("C10").Value = 0.002157599199563265
$ws.Range("D10").Value = 0.9998409152030945
$ws.Range("A11").Value = 0.0007988324505276978
$ws.Range("B11").Value = 0.9998267889022827
$ws.Range("C11").Value = 0.001679276814684272
$ws.Range("D11").Value = 0.9998409152030945
$ws.Range("A12").Value = 0.0005113819497637451
$ws.Range("B12").Value = 0.9998267889022827
$ws.Range("C12").Value = 0.002561234869062901
$ws.Range("D12").Value = 0.9998409152030945
$ws.Range("A13").Value = 0.0003231915470678359
$ws.Range("B13").Value = 0.9999133944511414
$ws.Range("C13").Value = 0.002869553165510297
$ws.Range("D13").Value = 0.9998409152030945
$ws.Range("A14").Value = 0.0002747249382082373
$ws.Range("B14").Value = 0.9999480247497559
$ws.Range("C14").Value = 0.002480077790096402
$ws.Range("D14").Value = 0.9998409152030945
$ws.Range("A15").Value = 0.0004984524566680193
$ws.Range("B15").Value = 0.9998614192008972
$ws.Range("C15").Value = 0.005236788187175989
$ws.Range("D15").Value = 0.9993637800216675
$ws.Range("A16").Value = 0.0007295574177987874
$ws.Range("B16").Value = 0.9998614192008972
$ws.Range("C16").Value = 0.003068211721256375
$ws.Range("D16").Value = 0.9996818900108337
$ws.Range("A17").Value = 0.0002776832261588424
$ws.Range("B17").Value = 0.9999133944511414
$ws.Range("C17").Value = 0.002934382064267993
$ws.Range("D17").Value = 0.9998409152030945
$ws.Range("A18").Value = 0.000207410441362299
$ws.Range("B18").Value = 0.999930739402771
$ws.Range("C18").Value = 0.002022790256887674
$ws.Range("D18").Value = 0.9998409152030945
$ws.Range("A19").Value = 0.00009976693399948999
$ws.Range("B19").Value = 0.9999480247497559
$ws.Range("C19").Value = 0.002919256454333663
$ws.Range("D19").Value = 0.9998409152030945
$ws.Range("A20").Value = 0.0007795181008987129
$ws.Range("B20").Value = 0.9998441338539124
$ws.Range("C20").Value = 0.007741168141365051
$ws.Range("D20").Value = 0.9993637800216675
$ws.Range("A21").Value = 0.000288108189124614
$ws.Range("B21").Value = 0.999930739402771
$ws.Range("C21").Value = 0.002888735383749008
$ws.Range("D21").Value = 0.9998409152030945
$ws.Range("A22").Value = 0.0001758837315719575
$ws.Range("B22").Value = 0.999930739402771
$ws.Range("C22").Value = 0.004064490552991629
$ws.Range("D22").Value = 0.9998409152030945
$ws.Range("A23").Value = 0.0001958083012141287
$ws.Range("B23").Value = 0.999930739402771
$ws.Range("C23").Value = 0.003381385933607817
$ws.Range("D23").Value = 0.9998409152030945
$ws.Range("A24").Value = 0.0001945407857419923
$ws.Range("B24").Value = 0.9999480247497559
$ws.Range("C24").Value = 0.001682174042798579
$ws.Range("D24").Value = 0.9998409152030945
$ws.Range("A25").Value = 0.00004539140354609117
$ws.Range("B25").Value = 0.9999826550483704
$ws.Range("C25").Value = 0.003241832600906491
$ws.Range("D25").Value = 0.9998409152030945
$ws.Range("A26").Value = 0.000003884140369336819
$ws.Range("B26").Value = 1
$ws.Range("C26").Value = 0.003300892654806376
$ws.Range("D26").Value = 0.9998409152030945
$ws.Range("A27").Value = 0.000387905485695228
$ws.Range("B27").Value = 0.9998787641525269
$ws.Range("C27").Value = 0.01055101212114096
$ws.Range("D27").Value = 0.9993637800216675
$ws.Range("A28").Value = 0.000413994857808575
$ws.Range("B28").Value = 0.9999480247497559
$ws.Range("C28").Value = 0.003726843278855085
$ws.Range("D28").Value = 0.9998409152030945
$ws.Range("A29").Value = 0.0002354870375711471
$ws.Range("B29").Value = 0.999930739402771
$ws.Range("C29").Value = 0.003831623122096062
$ws.Range("D29").Value = 0.9998409152030945
$ws.Range("A30").Value = 0.0005961621063761413
$ws.Range("B30").Value = 0.9998960494995117
$ws.Range("C30").Value = 0.002257796470075846
$ws.Range("D30").Value = 0.9998409152030945
$ws.Range("A31").Value = 0.0002154956018785015
$ws.Range("B31").Value = 0.9999653697013855
$ws.Range("C31").Value = 0.002923784544691443
$ws.Range("D31").Value = 0.9998409152030945
$ws.Range("A32").Value = 0.00007158840890042484
$ws.Range("B32").Value = 0.9999826550483704
$ws.Range("C32").Value = 0.00351674621924758
$ws.Range("D32").Value = 0.9998409152030945
$ws.Range("A33").Value = 0.00005663035699399188
$ws.Range("B33").Value = 0.9999826550483704
$ws.Range("C33").Value = 0.003559830598533154
$ws.Range("D33").Value = 0.9998409152030945
$ws.Range("A34").Value = 0.0007135042105801404
$ws.Range("B34").Value = 0.9998960494995117
$ws.Range("C34").Value = 0.003667050739750266
$ws.Range("D34").Value = 0.9998409152030945
$ws.Range("A35").Value = 0.00006492033571703359
$ws.Range("B35").Value = 0.9999653697013855
$ws.Range("C35").Value = 0.004719800315797329
$ws.Range("D35").Value = 0.9998409152030945
$ws.Range("A36").Value = 0.00004222943971399218
$ws.Range("B36").Value = 0.9999653697013855
$ws.Range("C36").Value = 0.005185364745557308
$ws.Range("D36").Value = 0.9998409152030945
$ws.Range("A37").Value = 0.0007182385306805372
$ws.Range("B37").Value = 0.9998960494995117
$ws.Range("C37").Value = 0.004706124775111675
$ws.Range("D37").Value = 0.9998409152030945
$ws.Range("A38").Value = 0.0002963297010865062
$ws.Range("B38").Value = 0.9999653697013855
$ws.Range("C38").Value = 0.002819443121552467
$ws.Range("D38").Value = 0.9998409152030945
$ws.Range("A39").Value = 0.00001656264612392988
$ws.Range("B39").Value = 1
$ws.Range("C39").Value = 0.003333253785967827
$ws.Range("D39").Value = 0.9998409152030945
$ws.Range("A40").Value = 0.00002228045559604652
$ws.Range("B40").Value = 0.9999826550483704
$ws.Range("C40").Value = 0.001885955804027617
$ws.Range("D40").Value = 0.9998409152030945
$ws.Range("A41").Value = 0.0001259826822206378
$ws.Range("B41").Value = 0.9999826550483704
$ws.Range("C41").Value = 0.004139212425798178
$ws.Range("D41").Value = 0.9998409152030945
$ws.Range("A42").Value = 0.0003364986623637378
$ws.Range("B42").Value = 0.9999480247497559
$ws.Range("C42").Value = 0.001209948444738984
$ws.Range("D42").Value = 0.9998409152030945
$ws.Range("A43").Value = 0.0004184081044513732
$ws.Range("B43").Value = 0.999930739402771
$ws.Range("C43").Value = 0.003793918527662754
$ws.Range("D43").Value = 0.9998409152030945
$ws.Range("A44").Value = 0.0001270915963687003
$ws.Range("B44").Value = 0.9999653697013855
$ws.Range("C44").Value = 0.002419582102447748
$ws.Range("D44").Value = 0.9998409152030945
$ws.Range("A45").Value = 0.00002880026113416534
$ws.Range("B45").Value = 0.9999826550483704
$ws.Range("C45").Value = 0.003060686634853482
$ws.Range("D45").Value = 0.9998409152030945
$ws.Range("A46").Value = 0.0002027652226388454
$ws.Range("B46").Value = 0.999930739402771
$ws.Range("C46").Value = 0.00438923854380846
$ws.Range("D46").Value = 0.9998409152030945
$ws.Range("A47").Value = 0.0000003863522408664721
$ws.Range("B47").Value = 1
$ws.Range("C47").Value = 0.004420455545186996
$ws.Range("D47").Value = 0.9998409152030945
$ws.Range("A48").Value = 0.0002706336381379515
$ws.Range("B48").Value = 0.9999653697013855
$ws.Range("C48").Value = 0.002682194812223315
$ws.Range("D48").Value = 0.9996818900108337
$ws.Range("A49").Value = 0.0001825051876949146
$ws.Range("B49").Value = 0.9999653697013855
$ws.Range("C49").Value = 0.003198240185156465
$ws.Range("D49").Value = 0.9998409152030945
$ws.Range("A50").Value = 0.00000732923899704474
$ws.Range("B50").Value = 1
$ws.Range("C50").Value = 0.003171787364408374
$ws.Range("D50").Value = 0.9998409152030945
$ws.Range("A51").Value = 0.0000004327342537635559
$ws.Range("C51").Value = 0.003686871146783233
$ws.Range("D51").Value = 0.9998409152030945
